$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Text/link columns (B, C) - plain string values, no numeric coercion risk
$ws.Range("B7").Value = "GateToken"
$ws.Range("C7").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("B8").Value = "MXToken"
$ws.Range("C8").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("B9").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C9").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("B10").Value = "WazirX"
$ws.Range("C10").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("B11").Value = "MCDex"
$ws.Range("C11").Value = "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
$ws.Range("B12").Value = "MandalaExchangeToken"
$ws.Range("C12").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("B13").Value = "BitrueCoin"
$ws.Range("C13").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("B14").Value = "BitMartToken"
$ws.Range("C14").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("B15").Value = "BitForexToken"
$ws.Range("C15").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("B16").Value = "CoinExToken"
$ws.Range("C16").Value = "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
$ws.Range("B17").Value = "TigerCash"
$ws.Range("C17").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("B18").Value = "LEO"
$ws.Range("C18").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"

# Numeric-looking text columns (D price, E percent) - force text format first
# so Excel does not coerce these into numbers, matching the original inline-string data.
$deCells = @("D2","E2","D3","E3","D4","E4","D5","E5","D6","E6","D7","E7","D8","E8","D9","E9","D10","E10","D11","E11","D12","E12","D13","E13","D14","E14","D15","E15","D16","E16","D17","E17","D18","E18","E19","D20","E20","D21","E21","D22","E22","E23","D24","E24","D25","E25","D26","E26","D38","E38","D39","E39","D40","E40","D41","E41","D42","E42","D43","E43","D44","E44","D45","E45","E46","D47","E47","D48","E48","D49","E49","D50","E50","D51","E51")
foreach ($addr in $deCells) { $ws.Range($addr).NumberFormat = "@" }

$ws.Range("D2").Value = "331.03"
$ws.Range("E2").Value = "0.54%"
$ws.Range("D3").Value = "44.19"
$ws.Range("E3").Value = "-0.11%"
$ws.Range("D4").Value = "5.539"
$ws.Range("E4").Value = "-0.58%"
$ws.Range("D5").Value = "0.08144"
$ws.Range("E5").Value = "0.54%"
$ws.Range("D6").Value = "2.058"
$ws.Range("E6").Value = "4.04%"
$ws.Range("D7").Value = "4.431"
$ws.Range("E7").Value = "2.40%"
$ws.Range("D8").Value = "0.9777"
$ws.Range("E8").Value = "2.66%"
$ws.Range("D9").Value = "0.1108"
$ws.Range("E9").Value = "-5.06%"
$ws.Range("D10").Value = "0.1893"
$ws.Range("E10").Value = "2.09%"
$ws.Range("D11").Value = "10.20"
$ws.Range("E11").Value = "-2.60%"
$ws.Range("D12").Value = "0.1001"
$ws.Range("E12").Value = "1.64%"
$ws.Range("D13").Value = "0.04719"
$ws.Range("E13").Value = "0.33%"
$ws.Range("D14").Value = "0.1058"
$ws.Range("E14").Value = "-0.90%"
$ws.Range("D15").Value = "0.001262"
$ws.Range("E15").Value = "-1.83%"
$ws.Range("D16").Value = "0.04098"
$ws.Range("E16").Value = "-3.07%"
$ws.Range("D17").Value = "0.005963"
$ws.Range("E17").Value = "1.02%"
$ws.Range("D18").Value = "3.343"
$ws.Range("E18").Value = "-0.91%"
$ws.Range("E19").Value = "1.69%"
$ws.Range("D20").Value = "0.3348"
$ws.Range("E20").Value = "-3.58%"
$ws.Range("D21").Value = "0.1389"
$ws.Range("E21").Value = "-1.43%"
$ws.Range("D22").Value = "0.2569"
$ws.Range("E22").Value = "2.46%"
$ws.Range("E23").Value = "3.86%"
$ws.Range("D24").Value = "0.004384"
$ws.Range("E24").Value = "0.74%"
$ws.Range("D25").Value = "0.0001278"
$ws.Range("E25").Value = "7.30%"
$ws.Range("D26").Value = "0.0003734"
$ws.Range("E26").Value = "-6.15%"
$ws.Range("D38").Value = "0.02689"
$ws.Range("E38").Value = "1.15%"
$ws.Range("D39").Value = "0.05647"
$ws.Range("E39").Value = "1.94%"
$ws.Range("D40").Value = "0.007630"
$ws.Range("E40").Value = "0.72%"
$ws.Range("D41").Value = "0.1421"
$ws.Range("E41").Value = "0.77%"
$ws.Range("D42").Value = "0.007561"
$ws.Range("E42").Value = "-6.45%"
$ws.Range("D43").Value = "0.001956"
$ws.Range("E43").Value = "-3.06%"
$ws.Range("D44").Value = "0.008296"
$ws.Range("E44").Value = "-6.87%"
$ws.Range("D45").Value = "0.00007015"
$ws.Range("E45").Value = "-3.37%"
$ws.Range("E46").Value = "-0.24%"
$ws.Range("D47").Value = "0.0005784"
$ws.Range("E47").Value = "-0.47%"
$ws.Range("D48").Value = "0.002516"
$ws.Range("E48").Value = "10.75%"
$ws.Range("D49").Value = "0.003525"
$ws.Range("E49").Value = "-26.03%"
$ws.Range("D50").Value = "0.00002097"
$ws.Range("E50").Value = "-0.24%"
$ws.Range("D51").Value = "0.0001997"
$ws.Range("E51").Value = "-0.24%"
